$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook previously computed each "bucket total" cell (B14:F20) with a
# long IF(year=2014, ..., IF(year=2015, ..., ...)) ladder that picked which
# merged year-header column to OFFSET from. That's replaced here with a
# single MATCH-against-row-1 lookup that finds the header column directly
# (and handles the last bucket, which has no "next" header, via IFERROR +
# COUNTA of row 2 to size the OFFSET).

$cols = @("B", "C", "D", "E", "F")

for ($i = 1; $i -le 7; $i++) {
    $row = 13 + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        if ($j + 1 -lt $cols.Length) { $nextCol = $cols[$j + 1] } else { $nextCol = "G" }
        $formula = '=SUM(OFFSET($A$1,ROWS($1:' + $i + '),MATCH(' + $col + '$13,$1:$1,0)-1,,IFERROR(MATCH(' + $nextCol + '$13,$1:$1,0),COUNTA($2:$2)+1)-MATCH(' + $col + '$13,$1:$1,0)))'
        $ws.Range($col + $row).Formula = $formula
    }
}

# New row 21: a totals row summing each bucket's column down through the
# newly-simplified rows 14:20.
$ws.Range("B21").Formula = "=SUM(B14:B20)"
$ws.Range("C21:F21").Formula = "=SUM(C14:C20)"
